# The presentation originally ends with a "Questions Time!" slide
# (slide 8, sldId 263 / rId9). This slide is being removed entirely,
# which is what the diff shows: slide8.xml (and its rels + content-type
# override) disappear, and the <p:sldId id="263" r:id="rId9"/> entry is
# dropped from the slide list.

$p = $ppt.ActivePresentation

# Locate the "Questions Time!" slide defensively (falls back to the
# last slide if, for some reason, the text can't be matched).
$targetIndex = $p.Slides.Count
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    foreach ($shape in $slide.Shapes) {
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            if ($shape.TextFrame.TextRange.Text -match "Questions Time!") {
                $targetIndex = $i
            }
        }
    }
}

$p.Slides.Item($targetIndex).Delete()
